# BVT checklist update for the "JSON Grid" visual.
# The "Basic chart" BVT row used to drive Values/Tooltip off 'Order Date' &
# 'Total Sales'; it is now driven off 'Date1' (tooltip) and
# 'Quarter'/'Month'/'Sales' (values), so rewrite the corresponding Steps /
# Output cells on the BVTs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")

# Row 3 ("Display default text" BVT) - Tooltip field now uses 'Date1'.
$ws.Range("D3").Value = "Drag 'Date1' column in 'Tooltip' field"

# Row 4 ("Display basic chart" BVT) - Values field now uses 'Quarter',
# 'Month' and 'Sales', and the rendered-visual confirmation text to match.
$ws.Range("D4").Value = "Drag 'Quarter', 'Month'  and 'Sales' column in 'Values' field"
$ws.Range("E4").Value = "Visual should be rendered for 'Quarter', 'Month'  and 'Sales'"

# Leave the sheet with D11 selected (matches the saved workbook's cursor).
$ws.Range("D11").Select()
